$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156, pushing existing rows 156:185 down to 157:186.
$ws.Rows("156:156").Insert()

# Populate the newly inserted row 156 with a new weekly price record,
# matching the surrounding data pattern (Femacal de La Calera / Coquimbo / Ciboulette).
$ws.Range("A156").Value = 3
$ws.Range("B156").Value = "Femacal de La Calera"
$ws.Range("C156").Value = "Coquimbo"
$ws.Range("D156").Value = 44476
$ws.Range("E156").Value = 5
$ws.Range("F156").Value = 100112039
$ws.Range("G156").Value = "Ciboulette"
$ws.Range("H156").Value = "Sin especificar"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 160
$ws.Range("K156").Value = 1500
$ws.Range("L156").Value = 1500
$ws.Range("M156").Value = 1500
$ws.Range("N156").Value = "$/docena de atados"
$ws.Range("O156").Value = "Provincia de Quillota"
$ws.Range("P156").Value = 500
$ws.Range("Q156").Value = 3
$ws.Range("R156").Value = "Hortaliza"
